$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 1-24 (columns A and B)
$data = @(
  @("b0003-1", "ACTTACGTGTCTGCGGTGTT"),
  @("b0003-2", "GGCCTCGTGATACGCCTATTGTCAGACTCCTAACTTCCATGAGA"),
  @("b0003-3", "CAAATAAAACGAAAGGCTCAGTACCCTCTCATGGAAGTTAGGAGTCTGACATGAAACTCTACAATCTGAAAGA"),
  @("b0003-4", "GCATTGTCTGTTGCGTGGTT"),
  @("b0003-7", "AATAGGCGTATCACGAGGCC"),
  @("b0003-8", "TGAGCCTTTCGTTTTATTTGATGC"),
  @("test-b0003-1", "TCAGATGCTGTTCAATACCGA"),
  @("test-b0003-2", "GTGCCGAGGAACAAGCCATA"),
  @("b0004-1", "CAATGACGTTACAGCTGCCG"),
  @("b0004-2", "GGCCTCGTGATACGCCTATTTTAGTTTTCCAGTACTCGTGCG"),
  @("b0004-3", "CAAATAAAACGAAAGGCTCAACGGCGGGCGCACGAGTACTGGAAAACTAAAATCTATTCATTATCTCAATCAGGC"),
  @("b0004-4", "ATCCGTCTTGAGAATGCCCG"),
  @("b0004-7", "AATAGGCGTATCACGAGGCC"),
  @("b0004-8", "TGAGCCTTTCGTTTTATTTGATGC"),
  @("test-b0004-1", "GCGTGAAGATTGCCGAAGTG"),
  @("test-b0004-2", "CGAAACCTTCAGCGAAGACG"),
  @("b0005-1", "TCTGTACATTGGGCGGCAAT"),
  @("b0005-2", "GGCCTCGTGATACGCCTATTAGGTTACTCCTTTTATCATTGCCC"),
  @("b0005-3", "CAAATAAAACGAAAGGCTCACATAACGGGCAATGATAAAAGGAGTAACCTATGACAAATGCCGGGTAACA"),
  @("b0005-4", "GACGCGTGGGATGATGTTTC"),
  @("b0005-7", "AATAGGCGTATCACGAGGCC"),
  @("b0005-8", "TGAGCCTTTCGTTTTATTTGATGC"),
  @("test-b0005-1", "TGACCATTCTGACCGCGAC"),
  @("test-b0005-2", "CGCCATTTTACTGCTCTCGC")
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $i + 1
  $ws.Cells.Item($row, 1).Value = $data[$i][0]
  $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove the now-unused rows 25-42
$ws.Range("A25:B42").EntireRow.Delete()
